# Review_345.docx: swap in the 11.11.24 review (Stealing Part of a
# Production Language Model, arXiv:2403.06634) in place of the
# 12.11.24 review (OccamLLM, arXiv:2406.06576). Each paragraphs
# run text is replaced in turn; the document structure (paragraph
# count, styles, runs) is left untouched.
$d = $word.ActiveDocument

# Paragraph 1: date in the daily-review header, 12.11.24 -> 11.11.24
$d.Content.Find.Execute('⚡️🚀המאמר היומי של מייק -12.11.24: ⚡️🚀', $false, $false, $false, $false, $false, $true, 1, $false, '⚡️🚀המאמר היומי של מייק -11.11.24: ⚡️🚀', 2) | Out-Null

# Paragraph 2: paper title
$d.Content.Find.Execute('OccamLLM: Fast and Exact Language Model Arithmetic in a Single Step', $false, $false, $false, $false, $false, $true, 1, $false, 'Stealing Part of a Production Language Model', 2) | Out-Null

# Paragraph 3: opening/intro paragraph of the review
$d.Content.Find.Execute('זהו מאמר שממש אהבתי, אהבתי גם את הרעיון וגם כתוב בצורה מאוד ברורה. למה כה אהבתי את הרעיון? אני כבר זמן מה טוען שבמקום להשקיע מאמצים גדולים באימון מודלי שפה לפתור בעיות מתמטיות יחסית מורכבות (שלדעתי מאוד קשה כי הם לא ״בנויים״ לזה באופן טבעי) כדאי להשתמש בכלים חיצוניים ייעודיים לכך (למשל כלים סימבוליים). מטרה של מודלי שפה במקרה הזה היא לזהות מתי הקלט שמוזן אליו (הפרומפט) מצריך פתרון בעיה מתמטית, ״לתרגם״ את הבעיה לשפה של הכלי הייעודי הזה, להעביר את הבעיה המתורגת לשפתו אליו לפתרון ולפענח את הפלט שלו.', $false, $false, $false, $false, $false, $true, 1, $false, 'מזמן לא סקרתי מאמר על איזה ניתן לפרוץ למודלים עמוקים. יש תחום שלם שנקרא adversarial learning שבו חוקרים מפתחים מנגנוני הגנה נגד התקפות שמנסות לגנוב משהו מהמודל או דרך המודל (למשל דאטה שהוא אומן עליו). המאמר שנסקור היום מציע שיטה שבאמצעותה ניתן לזהות המימד הפנימי (החבוי) של המודל (מימד ייצוגי הטוקנים) וגם את המטריצה בשכבה האחרונה של המודל. שכבה זו הממפה את האמבדינגס של כל הטוקנים ללוגיטים שלאחר מכן מוזנים לסופטמקס שממנו יוצרים ״ההסתברויות של הטוקנים.', 2) | Out-Null

# Paragraph 4: paragraph describing the paper's core idea
$d.Content.Find.Execute('וזה בדיוק מה שהמאמר הזה עושה. המחברים לקחו מודל שפה ופתחו מודל נפרד לפתרון בעיות מתמטיות. למעשה המודל לפתרון בעיות מתמטיות שפותח במאמר הוא גרף חישובי דינמי שכל צומת בו היא פונקציה או פעולה מתמטית (נדיג סימן + ו- *, או cos ו-exp). יש גם צמתים למשתני קלט השונים כדי שהמודל יוכל לחשב פונקציות על כמה משתנים (multivariate). למעשה גרף כזה הוא DAG או בשמו המלא Directed Acyclic Graph ומאמנים אותו לבחור את ״נתיב החישוב״ בו (״מסלול הצמתים״) בהינתן הייצוגים (אמבדינגס של הטוקנים) המוחשבים על ידי מודל שפה (ד״א מודל שפה לא מאומן ונותר קבוע לכל אורך אימון המודל). ', $false, $false, $false, $false, $false, $true, 1, $false, 'נתחיל מכך שמימד המטריצה W בשכבה האחרונה הוא N_voc x N_emb ,כאשר N_emb זה המימד הפנימי של המודל (אלפים בודדים) ו- N_voc הוא מספר הטוקנים במילון (בד״כ כמה עשרות אלפים ולפעמים מגיע מעל 100K).  כלומר N_voc > N_emb וזה בדיוק מה שמחברי המאמר מנצלים. מכיוון שהראנק של מטריצה W הוא N_emb כל המכפלות בה ממפות את הוקטור לתת מרחב במימד N_emb של מחרב הלוגיטים שהוא בעל מימד N_voc. כלומר אם ניקח מספר וקטורי לוגיטים ונשים אותם לעמודות של המטריצה (נקרא לה V) המספר המקסימלי לי וקטורים בלתי תלוים שיהיה לנו יהיה בדיוק N_emb.', 2) | Out-Null

# Paragraph 5: paragraph on how the authors exploit W's rank
$d.Content.Find.Execute('המחברים מאמנים שני מודלים: הראשון מזהה האם יש צורך בהפעלת המודל לחישובים מתמטיים לכל טוקן בהינתן ההקשר (כלומר כל הטוקנים לפניו). המודל השני מאומן לבנות נתיב חישובי בגרף החישובי שתיארתי בפסקה הקודמת. את שני המודלים האלו מאמנים בנפרד.', $false, $false, $false, $false, $false, $true, 1, $false, 'זה בדיוק מה שמחברי המאמר עשו. אולם מכיוון שהחישובים בטרנספורמרים הם לא בדיוק המלאה (FP16 גג) אז קשה לתפוס מתי העמודות הופכות להיות בלתי תלויות. במקום זה הם חישבו את הערכים הסינגולריים(ע״ס) של V (דרך מה שנקרא SVD - מי שלא מכיר ממליץ לקרוא על זה) ומסתכלים מתי היחס של ע״ס העוקבים (הם ממוינים) צונח משמעתית. ', 2) | Out-Null

# Paragraph 6: paragraph on detecting the hidden dimension
$d.Content.Find.Execute('מעניין כל שכבה של רשת ה-DAG הזה מורכבת משני חלקים: בחלק בראשון יש לנו צמתי החלטה: כל צומת כזה הוא וקטור ״המחבר״ אותו לצמתים פונקציונליים שכל אחד מהם הוא בעצם פעולה או פונקציה מתמטית (מקבוצת פעולות ופונקציות שבחרנו). הוקטור הזה הוא למעשה סופטמקס שממנו נדגם לאיזה צומת פונקציונלי/פעולה נחבר אותו. כל צומת פונקציונלי שנבחר מחובר עם כל צמתי ההחלטה מהשכבה הבאה ואליהם מועבר הייצוג משכבת ההחלטה הקודמת יחד עם ייצוג הפעולה (כנראה האם נבחרה או לא). כך נבנה גרף חישובי מייצוגי הטוקנים המחושבים על ידי מודל שפה (הם מחוברים לשכבת ההחלטה הראשון במודל החישובי). ד״א כל פעולה וכל פונקציית בסיס בגרף משוכפלת בכמה צמתית כדי להקנות למודל יכולת לקרב פונקציות מורכבות יותר.', $false, $false, $false, $false, $false, $true, 1, $false, 'למה זה חשוב? כי במקרה האידאלי ע״ס של V צריכים להתאפס אחר שעברנו את הראנק של או N_emb. אז בגלל אי דיוקים נומריים במודל כמובן שלא נראה ממש אפסים שם אלא ערכים מאוד נמוכים ואיפה שזה מתחיל לקרות זה בדיוק במימד N_emb + 1. אז עושים את הטריק הזה על הרבה מאוד דאטה ומגלים את המימד החבוי של המודל שלכם.', 2) | Out-Null

# Paragraph 7: paragraph on the topK logits trick (contains an apostrophe,
# so it is set directly via Range.Text -- Find.Execute's ReplaceWith
# smart-quotes a plain apostrophe into a curly one)
$d.Paragraphs(7).Range.Text = 'כמובן שבעולם האמיתי אין לכם גישה לכל הלוגיטים אלא רק ל-topK ואז המאמר מנצל את העובדה שניתן לקנפג חלק מהמודל להוסיף מרג''ין לטוקן נתון במילון. ואחרי מספיק משחקים מקבלים את כל הלוגיטים (זה די יקר חישובית).'

# Paragraph 8: closing paragraph / verdict
$d.Content.Find.Execute('מאמר די נחמד אבל כתוב לא מאוד ברור (או שהיה חסר לי קצת רקע)...', $false, $false, $false, $false, $false, $true, 1, $false, 'מימד של W זה נחמד אבל מה עם מטריצה W עצמה. המאמר מציע התקפה כדי לגלות אותה (סוג של) גם. בכללי המאמר מלא ברעיונות יפים להתקפות על המודלים ומי שמתעניין מוזמן להעיף מבט.', 2) | Out-Null

# Paragraph 9: arXiv link
$d.Content.Find.Execute('https://arxiv.org/abs/2406.06576', $false, $false, $false, $false, $false, $true, 1, $false, 'https://arxiv.org/abs/2403.06634', 2) | Out-Null
